# Commit: theme swap ("Integral" -> default "Office Theme" colors on the
# live slide-master theme part) plus a table style-id change on the
# sources-of-finance table on slide 6.
#
# Converts a hex colour string ("RRGGBB") to the little-endian Long that
# VBA's RGB()/ColorFormat.RGB expects (0xBBGGRR).
function HexToRgbInt([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Table on slide 6 ("SOURCES OF FINANCE ...") switches to a different
#    built-in table style GUID.
# ---------------------------------------------------------------------
$slide6 = $p.Slides.Item(6)
$tableShape = $slide6.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{9CA4FB1C-381B-48C1-93C4-E42746AF359A}")

# ---------------------------------------------------------------------
# 2) The presentation's working theme (the one behind the slide master)
#    swaps from the "Integral" palette to the stock "Office" palette.
#    MsoThemeColorSchemeIndex 1-12 maps to dk1, lt1, dk2, lt2, accent1-6,
#    hlink, folHlink - exactly the order/values of the "Office" scheme.
# ---------------------------------------------------------------------
$officeColors = @(
    "000000", # 1  dk1
    "FFFFFF", # 2  lt1
    "44546A", # 3  dk2
    "E7E6E6", # 4  lt2
    "5B9BD5", # 5  accent1
    "ED7D31", # 6  accent2
    "A5A5A5", # 7  accent3
    "FFC000", # 8  accent4
    "4472C4", # 9  accent5
    "70AD47", # 10 accent6
    "0563C1", # 11 hlink
    "954F72"  # 12 folHlink
)

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $officeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = HexToRgbInt($officeColors[$i - 1])
}
